$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 61.257969
$ws.Range("H2").Value = 122.515938
$ws.Range("I2").Value = 0.1582655541224298
$ws.Range("J2").Value = 0.1141811124301534
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.032708
$ws.Range("N2").Value = 0.065416
$ws.Range("Q2").Value = 2.003625650052
$ws.Range("R2").Value = 8.014502600208001
$ws.Range("S2").Value = 0.1582655541224298
$ws.Range("T2").Value = 0.1141811124301534

# Row 3
$ws.Range("I3").Value = 0.3520419132016297
$ws.Range("J3").Value = 0.3809723868306807
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.032708
$ws.Range("N3").Value = 0.065416
$ws.Range("Q3").Value = 4.456814441369334
$ws.Range("R3").Value = 26.740886648216
$ws.Range("S3").Value = 0.3520419132016297
$ws.Range("T3").Value = 0.3809723868306807

# Row 4
$ws.Range("G4").Value = 56.53322600000001
$ws.Range("H4").Value = 169.599678
$ws.Range("I4").Value = 0.1460587493395114
$ws.Range("J4").Value = 0.1580617201154336
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.032708
$ws.Range("N4").Value = 0.065416
$ws.Range("Q4").Value = 1.849088756008
$ws.Range("R4").Value = 11.094532536048
$ws.Range("S4").Value = 0.1460587493395114
$ws.Range("T4").Value = 0.1580617201154336

# Row 5
$ws.Range("G5").Value = 26.919878
$ws.Range("H5").Value = 53.83975599999999
$ws.Range("I5").Value = 0.06954996187644104
$ws.Range("J5").Value = 0.05017700826032956
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.032708
$ws.Range("N5").Value = 0.065416
$ws.Range("Q5").Value = 0.8804953696239999
$ws.Range("R5").Value = 3.521981478496
$ws.Range("S5").Value = 0.06954996187644104
$ws.Range("T5").Value = 0.05017700826032956

# Row 6
$ws.Range("G6").Value = 47.55688233333333
$ws.Range("H6").Value = 142.670647
$ws.Range("I6").Value = 0.1228675461770565
$ws.Range("J6").Value = 0.1329646856688125
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.032708
$ws.Range("N6").Value = 0.065416
$ws.Range("Q6").Value = 1.555490507358667
$ws.Range("R6").Value = 9.332943044152
$ws.Range("S6").Value = 0.1228675461770565
$ws.Range("T6").Value = 0.1329646856688125

# Row 7
$ws.Range("G7").Value = 58.52948833333333
$ws.Range("H7").Value = 175.588465
$ws.Range("I7").Value = 0.1512162752829316
$ws.Range("J7").Value = 0.1636430866945903
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.032708
$ws.Range("N7").Value = 0.065416
$ws.Range("Q7").Value = 1.914382504406666
$ws.Range("R7").Value = 11.48629502644
$ws.Range("S7").Value = 0.1512162752829316
$ws.Range("T7").Value = 0.1636430866945903
